# Update the build version string throughout the workbook.
# Old version string: "mines - January 30 (built on January 30 2026 16.19.47 EST)"
# New version string: "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: $newVersion"

$aboutWs.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Oak Grove Mine, United States, M3577, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

$usedRange = $dataWs.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataWs.Cells.Item($r, 19)  # Column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
